$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 600
$ws.Range("J7").Value = 600
$ws.Range("L7").Value = 600
$ws.Range("N7").Value = -824

# Row 14
$ws.Range("H14").Value = 600
$ws.Range("J14").Value = 600
$ws.Range("L14").Value = 600
$ws.Range("N14").Value = -982

# Row 17
$ws.Range("H17").Value = 1037782.06
$ws.Range("J17").Value = 1578990.9
$ws.Range("L17").Value = 4736972.699999999
$ws.Range("N17").Value = -4737308.699999999

# Row 74
$ws.Range("H74").Value = 4143.8887
$ws.Range("I74").Value = 4080.5881
$ws.Range("J74").Value = 4251.5
$ws.Range("K74").Value = 4080.5881
$ws.Range("L74").Value = 4251.5
$ws.Range("M74").Value = -3144.5881
$ws.Range("N74").Value = -6123.5

# Row 77
$ws.Range("H77").Value = 4143.8887
$ws.Range("I77").Value = 4080.5881
$ws.Range("J77").Value = 4251.5
$ws.Range("K77").Value = 20402.9405
$ws.Range("L77").Value = 21257.5
$ws.Range("M77").Value = -15722.9405
$ws.Range("N77").Value = -30617.5

# Row 95
$ws.Range("H95").Value = 50333.332
$ws.Range("J95").Value = 50333.332
$ws.Range("L95").Value = 50333.332
$ws.Range("N95").Value = -55825.332

# Row 112
$ws.Range("H112").Value = 5389.0293
$ws.Range("I112").Value = 375
$ws.Range("K112").Value = 1125
$ws.Range("M112").Value = -17

# Row 113
$ws.Range("H113").Value = 3102
$ws.Range("I113").Value = 3461
$ws.Range("J113").Value = 2877.625
$ws.Range("K113").Value = 3461
$ws.Range("L113").Value = 2877.625
$ws.Range("M113").Value = -207
$ws.Range("N113").Value = -9385.625

# Row 115
$ws.Range("H115").Value = 2326.6191

# Row 138
$ws.Range("H138").Value = 2661.8076
$ws.Range("I138").Value = 2305.5715
$ws.Range("J138").Value = 4158
$ws.Range("K138").Value = 6916.7145
$ws.Range("L138").Value = 12474
$ws.Range("M138").Value = -1776.7145
$ws.Range("N138").Value = -22754


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1127091.2
$ws.Range("I32").Value = 1326552.9
$ws.Range("J32").Value = 30052.334
$ws.Range("K32").Value = 1326552.9
$ws.Range("L32").Value = 30052.334
$ws.Range("M32").Value = -1326265.9
$ws.Range("N32").Value = -30626.334

# Row 88
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()

# Row 91
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()


$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4980
$ws.Range("I86").Value = 4980
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4980
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3857
$ws.Range("N86").ClearContents()

# Row 89
$ws.Range("H89").Value = 4980
$ws.Range("I89").Value = 4980
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 24900
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -19284
$ws.Range("N89").ClearContents()

# Row 107
$ws.Range("H107").Value = 63651.062
$ws.Range("I107").Value = 84293.086
$ws.Range("J107").Value = 1725
$ws.Range("K107").Value = 84293.086
$ws.Range("L107").Value = 1725
$ws.Range("M107").Value = -82373.086
$ws.Range("N107").Value = -5565


$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 39500.332
$ws.Range("J4").Value = 39500.332
$ws.Range("L4").Value = 39500.332
$ws.Range("N4").Value = -39724.332

# Row 107
$ws.Range("H107").Value = 3290234
$ws.Range("I107").Value = 10416946
$ws.Range("K107").Value = 10416946
$ws.Range("M107").Value = -10415026


$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 24000836
$ws.Range("I4").Value = 17500796
$ws.Range("J4").Value = 50001000
$ws.Range("K4").Value = 52502388
$ws.Range("L4").Value = 150003000
$ws.Range("M4").Value = -52502276
$ws.Range("N4").Value = -150003224

# Row 63
$ws.Range("H63").Value = 3485.3333

# Row 66
$ws.Range("H66").Value = 3485.3333

# Row 137
$ws.Range("H137").Value = 20849434
$ws.Range("I137").Value = 83383336
$ws.Range("J137").Value = 4800
$ws.Range("K137").Value = 250150008
$ws.Range("L137").Value = 14400
$ws.Range("M137").Value = -250144908
$ws.Range("N137").Value = -24600

# Row 140
$ws.Range("H140").Value = 1344.8536
$ws.Range("I140").Value = 943.4138
$ws.Range("J140").Value = 2315
$ws.Range("K140").Value = 2830.2414
$ws.Range("L140").Value = 6945
$ws.Range("M140").Value = 2349.7586
$ws.Range("N140").Value = -17305


$ws = $wb.Worksheets.Item("GSM")
# Row 47
$ws.Range("H47").Value = 5000
$ws.Range("J47").Value = 5000
$ws.Range("L47").Value = 5000
$ws.Range("N47").Value = -6136

# Row 80
$ws.Range("H80").Value = 1402068.6
$ws.Range("I80").Value = 4502250
$ws.Range("J80").Value = 161996
$ws.Range("K80").Value = 4502250
$ws.Range("L80").Value = 161996
$ws.Range("M80").Value = -4501252
$ws.Range("N80").Value = -163992

# Row 82
$ws.Range("H82").Value = 24955.562
$ws.Range("J82").Value = 24999.4
$ws.Range("L82").Value = 24999.4
$ws.Range("N82").Value = -25765.4

# Row 83
$ws.Range("H83").Value = 1402068.6
$ws.Range("I83").Value = 4502250
$ws.Range("J83").Value = 161996
$ws.Range("K83").Value = 22511250
$ws.Range("L83").Value = 809980
$ws.Range("M83").Value = -22506258
$ws.Range("N83").Value = -819964

# Row 85
$ws.Range("H85").Value = 24955.562
$ws.Range("J85").Value = 24999.4
$ws.Range("L85").Value = 24999.4
$ws.Range("N85").Value = -27651.4

# Row 97
$ws.Range("H97").Value = 522.6
$ws.Range("I97").Value = 502.7857
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 502.7857
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -6.78570000000002
$ws.Range("N97").Value = -1792


$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 23519.334
$ws.Range("I62").Value = 3019.8
$ws.Range("J62").Value = 33769.1
$ws.Range("K62").Value = 3019.8
$ws.Range("L62").Value = 33769.1
$ws.Range("M62").Value = -2395.8
$ws.Range("N62").Value = -35017.1

# Row 65
$ws.Range("H65").Value = 23519.334
$ws.Range("I65").Value = 3019.8
$ws.Range("J65").Value = 33769.1
$ws.Range("K65").Value = 15099
$ws.Range("L65").Value = 168845.5
$ws.Range("M65").Value = -11979
$ws.Range("N65").Value = -175085.5

# Row 81
$ws.Range("H81").Value = 3451.1052
$ws.Range("J81").Value = 2766.7778
$ws.Range("L81").Value = 5533.5556
$ws.Range("N81").Value = -7655.5556

# Row 84
$ws.Range("H84").Value = 3451.1052
$ws.Range("J84").Value = 2766.7778
$ws.Range("L84").Value = 27667.778
$ws.Range("N84").Value = -38275.778

# Row 132
$ws.Range("H132").Value = 3146783.5
$ws.Range("I132").Value = 2615
$ws.Range("J132").Value = 5954077
$ws.Range("K132").Value = 7845
$ws.Range("L132").Value = 17862231
$ws.Range("M132").Value = -5315
$ws.Range("N132").Value = -17867291

# Row 136
$ws.Range("H136").Value = 1957.5916
$ws.Range("I136").Value = 1598.5178
$ws.Range("K136").Value = 4795.553400000001
$ws.Range("M136").Value = -2245.553400000001

